$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.940.25'
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").Value = '1.984.53'
$ws.Range("E3").Value = '  +1.01%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'245.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.628"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.62%  '
$ws.Range("D7").Value = "'60.84"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.73%  '
$ws.Range("E9").Value = '  +2.03%  '
$ws.Range("D10").Value = "'0.0799"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.86%  '
$ws.Range("D11").Value = "'0.103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("E12").Value = '  +9.15%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = "'0.845"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.90%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = "'22.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("D15").Value = '2.276.94'
$ws.Range("E15").Value = '  +1.08%  '
$ws.Range("D16").Value = "'5.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.65%  '
$ws.Range("D17").Value = '1.983.07'
$ws.Range("E17").Value = '  +0.81%  '
$ws.Range("D18").Value = '36.849.70'
$ws.Range("E18").Value = '  +0.85%  '
$ws.Range("D19").Value = "'70.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("D20").Value = '0.0₃0861'
$ws.Range("E20").Value = '  +0.22%  '
$ws.Range("D21").Value = "'5.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.14%  '
$ws.Range("D22").Value = "'229.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  +2.66%  '
$ws.Range("D25").Value = "'2.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.61%  '
$ws.Range("D26").Value = "'0.149"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.33%  '
$ws.Range("D27").Value = "'9.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.76%  '
$ws.Range("D28").Value = "'163.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.77%  '
$ws.Range("E29").Value = '  +0.59%  '
$ws.Range("D30").Value = "'1.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +20.02%  '
$ws.Range("E31").Value = '  +1.72%  '
$ws.Range("D32").Value = "'4.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.43%  '
$ws.Range("D33").Value = "'0.0621"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("D34").Value = "'4.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.02%  '
$ws.Range("D35").Value = "'2.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.93%  '
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").Value = "'3.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("E38").Value = '  +0.19%  '
$ws.Range("D39").Value = "'5.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.08%  '
$ws.Range("D40").Value = "'0.0996"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.82%  '
$ws.Range("E41").Value = '  +0.75%  '
$ws.Range("E42").Value = '  +0.61%  '
$ws.Range("E43").Value = '  +0.60%  '
$ws.Range("D44").Value = "'16.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.73%  '
$ws.Range("D45").Value = "'90.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.64%  '
$ws.Range("D46").Value = '1.368.46'
$ws.Range("E46").Value = '  +0.51%  '
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").Value = "'7.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D49").Value = "'46.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.11%  '
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").Value = "'2.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.78%  '
$ws.Range("E51").Value = '  +11.03%  '
